# "Inclusao do driver chrome para webdriver"
#
# The "Convenio" (agreement) test scenarios are being dropped from the test
# script in favor of an "excluirEmpresa" (delete company) scenario reusing
# row 4, which previously exercised "abreFormConvenio"/"formConvenio".
#
# Net effect on the sheet:
#   - rows 5-9 (cadastrarConvenio / excluirConvenio / old abreFormEmpresa /
#     excluirEmpresa test rows) are removed entirely
#   - row 4 becomes the new "excluirEmpresa" test case: action + landing
#     state columns are rewritten and a CNPJ id is added in column E,
#     while the now-unused trailing columns F:N are cleared out

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old convenio-related rows (5 through 9) completely - this shifts
# nothing below them since they are the last rows, and it naturally drops
# the shared strings that only they referenced.
$ws.Rows("5:9").Delete()

# Row 4 used to be "abreFormConvenio" / "formConvenio" with no data columns
# populated; it becomes "excluirEmpresa" / "formEmpresa" with a CNPJ value
# in column E. Clear out the old trailing blank/styled cells first.
$ws.Range("F4:N4").Clear()
$ws.Range("C4:D4").ClearFormats()

$ws.Range("C4").Value2 = "excluirEmpresa"
$ws.Range("D4").Value2 = "formEmpresa"

# Write the CNPJ as text (leading apostrophe forces text, matching the
# quote-prefixed numeric-id style used elsewhere in the sheet, e.g. F3),
# then pick up that same cell style so it matches the rest of the sheet.
$ws.Range("E4").Value = "'89424232000180"
$ws.Range("F3").Copy()
$ws.Range("E4").PasteSpecial(-4122)
$excel.CutCopyMode = 0
